# Refresh the cryptos price/volume table with the latest scrape.
# Some "Price" (column D) values are plain decimal numbers (e.g. "0.537");
# assigning those straight to .Value would let Excel auto-convert them to a
# numeric type and drop the original text formatting (e.g. trailing zeros).
# Set-TextValue forces the cell's number format to Text ("@") just for the
# write, then restores the cell's original Style so no formatting drifts.
function Set-TextValue {
    param($ws, $addr, $val)
    $origStyle = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.396.59"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "2.428.75"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  +0.39%  "
Set-TextValue $ws "D5" "556.99"
$ws.Range("E5").Value = "  +2.18%  "
Set-TextValue $ws "D6" "143.91"
$ws.Range("E6").Value = "  +5.27%  "
$ws.Range("E7").Value = "  +0.09%  "
Set-TextValue $ws "D8" "0.537"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").Value = "2.429.73"
$ws.Range("E10").Value = "  +5.16%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  +1.64%  "
Set-TextValue $ws "D13" "0.352"
$ws.Range("E13").Value = "  +3.74%  "
Set-TextValue $ws "D14" "26.33"
$ws.Range("E14").Value = "  +6.91%  "
$ws.Range("E15").Value = "  +9.70%  "
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "62.179.91"
$ws.Range("E17").Value = "  +2.45%  "
$ws.Range("D18").Value = "2.428.16"
$ws.Range("E18").Value = "  +3.51%  "
Set-TextValue $ws "D19" "11.11"
$ws.Range("E19").Value = "  +4.67%  "
Set-TextValue $ws "D20" "325.29"
$ws.Range("E20").Value = "  +1.91%  "
Set-TextValue $ws "D21" "4.17"
$ws.Range("E21").Value = "  +1.67%  "
Set-TextValue $ws "D22" "6.77"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("E23").Value = "  +1.15%  "
Set-TextValue $ws "D24" "1.79"
$ws.Range("E24").Value = "  +3.12%  "
Set-TextValue $ws "D25" "65.07"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("E26").Value = "  +10.38%  "
Set-TextValue $ws "D27" "573.02"
$ws.Range("E27").Value = "  +15.10%  "
$ws.Range("D28").Value = "2.538.26"
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0948"
$ws.Range("E29").Value = "  +9.44%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws "D30" "0.999"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +5.87%  "
$ws.Range("E32").Value = "  +5.78%  "
$ws.Range("E33").Value = "  +2.30%  "
Set-TextValue $ws "D34" "1.86"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("E35").Value = "  +4.74%  "
Set-TextValue $ws "D36" "5.74"
$ws.Range("E36").Value = "  +8.90%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +5.31%  "
$ws.Range("E39").Value = "  +2.55%  "
Set-TextValue $ws "D40" "1.89"
$ws.Range("E40").Value = "  +4.42%  "
Set-TextValue $ws "D41" "18.88"
$ws.Range("E41").Value = "  +1.74%  "
Set-TextValue $ws "D42" "147.68"
$ws.Range("E42").Value = "  +4.45%  "
$ws.Range("E43").Value = "  +0.42%  "
Set-TextValue $ws "D44" "41.68"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("E45").Value = "  +12.92%  "
Set-TextValue $ws "D46" "152.14"
$ws.Range("E46").Value = "  +6.89%  "
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("E48").Value = "  +6.00%  "
Set-TextValue $ws "D49" "20.49"
$ws.Range("E49").Value = "  +7.57%  "
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("E51").Value = "  +4.41%  "
